$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4982.778
$ws.Range("I113").Value = 4140.9
$ws.Range("J113").Value = 6035.125
$ws.Range("K113").Value = 4140.9
$ws.Range("L113").Value = 6035.125
$ws.Range("M113").Value = -886.8999999999996
$ws.Range("N113").Value = -12543.125
$ws.Range("H132").Value = 2157.5833
$ws.Range("I132").Value = 1518.4193
$ws.Range("J132").Value = 6120.4
$ws.Range("K132").Value = 4555.257900000001
$ws.Range("L132").Value = 18361.2
$ws.Range("M132").Value = -2025.257900000001
$ws.Range("N132").Value = -23421.2
$ws.Range("H135").Value = 2115.7646
$ws.Range("I135").Value = 2141.5334
$ws.Range("K135").Value = 19273.8006
$ws.Range("M135").Value = -16738.8006

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 64024
$ws.Range("J10").Value = 64024
$ws.Range("L10").Value = 64024
$ws.Range("N10").Value = -64364
$ws.Range("H14").Value = 701
$ws.Range("I14").Value = 701
$ws.Range("K14").Value = 701
$ws.Range("M14").Value = -526
$ws.Range("H16").Value = 13999.333
$ws.Range("I16").Value = 14998
$ws.Range("J16").Value = 13500
$ws.Range("K16").Value = 14998
$ws.Range("L16").Value = 13500
$ws.Range("M16").Value = -14711
$ws.Range("N16").Value = -14074
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H18").Value = 7210
$ws.Range("J18").Value = 7210
$ws.Range("L18").Value = 7210
$ws.Range("N18").Value = -7854
$ws.Range("H19").Value = 2632
$ws.Range("I19").Value = 1448
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 1448
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -1219
$ws.Range("N19").Value = -5458
$ws.Range("H21").Value = 4609.75
$ws.Range("I21").Value = 1999.5
$ws.Range("J21").Value = 7220
$ws.Range("K21").Value = 1999.5
$ws.Range("L21").Value = 7220
$ws.Range("M21").Value = -1625.5
$ws.Range("N21").Value = -7968
$ws.Range("H25").Value = 1034.8572
$ws.Range("I25").Value = 1034.8572
$ws.Range("K25").Value = 1034.8572
$ws.Range("M25").Value = -632.8571999999999
$ws.Range("H27").Value = 8210
$ws.Range("J27").Value = 8210
$ws.Range("L27").Value = 8210
$ws.Range("N27").Value = -8578
$ws.Range("H30").Value = 4864
$ws.Range("I30").Value = 1350
$ws.Range("J30").Value = 7206.6665
$ws.Range("K30").Value = 1350
$ws.Range("L30").Value = 7206.6665
$ws.Range("M30").Value = -1200
$ws.Range("N30").Value = -7506.6665
$ws.Range("H61").Value = 3116.2
$ws.Range("I61").Value = 3009.5386
$ws.Range("J61").Value = 3314.2856
$ws.Range("K61").Value = 3009.5386
$ws.Range("L61").Value = 3314.2856
$ws.Range("M61").Value = -2797.5386
$ws.Range("N61").Value = -3738.2856
$ws.Range("H136").Value = 3116.2
$ws.Range("I136").Value = 3009.5386
$ws.Range("J136").Value = 3314.2856
$ws.Range("K136").Value = 9028.6158
$ws.Range("L136").Value = 9942.856800000001
$ws.Range("M136").Value = -6478.6158
$ws.Range("N136").Value = -15042.8568

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 14444.333
$ws.Range("I18").Value = 3000
$ws.Range("J18").Value = 20166.5
$ws.Range("K18").Value = 3000
$ws.Range("L18").Value = 20166.5
$ws.Range("M18").Value = -2471
$ws.Range("N18").Value = -21224.5
$ws.Range("H20").Value = 8390.638999999999
$ws.Range("I20").Value = 1282.7084
$ws.Range("J20").Value = 22606.5
$ws.Range("K20").Value = 1282.7084
$ws.Range("L20").Value = 22606.5
$ws.Range("M20").Value = -1035.7084
$ws.Range("N20").Value = -23100.5
$ws.Range("H23").Value = 3986
$ws.Range("I23").Value = 762
$ws.Range("J23").Value = 7210
$ws.Range("K23").Value = 762
$ws.Range("L23").Value = 7210
$ws.Range("M23").Value = -479
$ws.Range("N23").Value = -7776

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -50
$ws.Range("H132").Value = 7496
$ws.Range("I132").Value = 7496
$ws.Range("K132").Value = 22488
$ws.Range("M132").Value = -19958

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 55
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 165
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -52
$ws.Range("N6").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 174.16667
$ws.Range("I2").Value = 198.5
$ws.Range("J2").Value = 125.5
$ws.Range("K2").Value = 198.5
$ws.Range("L2").Value = 125.5
$ws.Range("M2").Value = -85.5
$ws.Range("N2").Value = -351.5
$ws.Range("H18").Value = 3150
$ws.Range("I18").Value = 3150
$ws.Range("K18").Value = 3150
$ws.Range("M18").Value = -2857
$ws.Range("H21").Value = 40000
$ws.Range("J21").Value = 40000
$ws.Range("L21").Value = 40000
$ws.Range("N21").Value = -40346
$ws.Range("H30").Value = 40000
$ws.Range("J30").Value = 40000
$ws.Range("L30").Value = 40000
$ws.Range("N30").Value = -40210
$ws.Range("H58").Value = 11800
$ws.Range("I58").Value = 9700
$ws.Range("J58").Value = 16000
$ws.Range("K58").Value = 9700
$ws.Range("L58").Value = 16000
$ws.Range("M58").Value = -9423
$ws.Range("N58").Value = -16554

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1009
$ws.Range("I4").Value = 1009
$ws.Range("K4").Value = 1009
$ws.Range("M4").Value = -896
$ws.Range("H28").Value = 1009
$ws.Range("I28").Value = 1009
$ws.Range("K28").Value = 1009
$ws.Range("M28").Value = -777
$ws.Range("H37").Value = 1009
$ws.Range("I37").Value = 1009
$ws.Range("K37").Value = 1009
$ws.Range("M37").Value = -902
$ws.Range("H136").Value = 6480.1035
$ws.Range("I136").Value = 3191.7727
$ws.Range("J136").Value = 16814.857
$ws.Range("K136").Value = 9575.3181
$ws.Range("L136").Value = 50444.571
$ws.Range("M136").Value = -7025.3181
$ws.Range("N136").Value = -55544.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1160
$ws.Range("I136").Value = 1036.3636
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 3109.0908
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -559.0907999999999
$ws.Range("N136").Value = -9600
